$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 302; all rows from 302-392 shift down to 303-393,
# and the sheet's used range grows from R392 to R393.
$ws.Rows("302:302").Insert()

# Populate the newly inserted row 302 with the new data record.
$ws.Cells.Item(302, 1).Value  = 9
$ws.Cells.Item(302, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(302, 3).Value  = "Metropolitana"
$ws.Cells.Item(302, 4).Value  = 44736
$ws.Cells.Item(302, 5).Value  = 13
$ws.Cells.Item(302, 6).Value  = 100112044
$ws.Cells.Item(302, 7).Value  = "Perejil"
$ws.Cells.Item(302, 8).Value  = "Sin especificar"
$ws.Cells.Item(302, 9).Value  = "Primera"
$ws.Cells.Item(302, 10).Value = 52
$ws.Cells.Item(302, 11).Value = 12000
$ws.Cells.Item(302, 12).Value = 13000
$ws.Cells.Item(302, 13).Value = 12500
$ws.Cells.Item(302, 14).Value = "`$/docena de atados"
$ws.Cells.Item(302, 15).Value = "Región Metropolitana"
$ws.Cells.Item(302, 16).Value = 4167
$ws.Cells.Item(302, 17).Value = 3
$ws.Cells.Item(302, 18).Value = "Hortaliza"
